$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames / new column -----------------------------------------
# Order matters for shared-string table ordering: Notes, then
# "Sensor EquipmentCode", then "DistanceFromProbeTip_cm".
$ws.Range("H7").Value2 = "Notes"
$ws.Range("E7").Value2 = "Sensor EquipmentCode"
$ws.Range("F7").Value2 = "DistanceFromProbeTip_cm"

# Header row is now taller (wrapped two-line captions).
$ws.Range("A7").RowHeight = 39

# Column E (Sensor EquipmentCode) is narrower now.
$ws.Columns.Item(5).ColumnWidth = 7.71

# --- Data rows: shift DistanceFromProbeTip_cm by +2 cm and compute the ---
# --- new Depth_cm column (G) = 130 - 28 - DistanceFromProbeTip_cm --------
for ($r = 8; $r -le 43; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $fCell.Value2 = $fCell.Value2 + 2

    $gCell = $ws.Cells.Item($r, 7)
    $gCell.Formula = "=130-28-F$r"
    $gCell.Font.Color = 16711680
}

# --- Freeze panes below the header / right of the SiteCode columns -------
$ws.Range("D8").Select()
$excel.ActiveWindow.FreezePanes = $true

# Leave each pane with its own remembered selection, as the authored sheet
# had: top-right pane on D1, bottom-left pane on A8, bottom-right (the
# active pane) on I28.
$ws.Range("D1").Select()
$ws.Range("A8").Select()
$ws.Range("I28").Select()
